$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.473.93'
$ws.Range('E2').Value = '  +2.51%  '

$ws.Range('D3').Value = '2.313.30'
$ws.Range('E3').Value = '  +1.69%  '

$ws.Range('E4').Value = '  -0.05%  '

$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '311.37'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +1.64%  '

$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '102.50'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +4.65%  '

$ws.Range('E7').Value = '  +1.50%  '

$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.532'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  +7.27%  '

$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '35.74'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  +0.92%  '

$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.0817'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  +3.03%  '

$ws.Range('E12').Value = '  -0.45%  '

$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '7.02'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  +1.61%  '

$ws.Range('D14').Value = '2.672.97'
$ws.Range('E14').Value = '  +1.65%  '

$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '15.02'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  +1.43%  '

$ws.Range('D16').Value = '2.319.81'
$ws.Range('E16').Value = '  +1.88%  '

$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '0.810'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  +1.66%  '

$ws.Range('D18').Value = '43.381.62'
$ws.Range('E18').Value = '  +2.63%  '

$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '12.43'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  -0.76%  '

$ws.Range('D20').Value = '0.0₃0926'
$ws.Range('E20').Value = '  +2.00%  '

$ws.Range('E21').Value = '  +1.93%  '

$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '68.22'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -0.08%  '

$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '241.75'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +1.14%  '

$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '2.03'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +3.29%  '

$ws.Range('E25').Value = '  +1.77%  '

$ws.Range('E26').Value = '  -0.08%  '

$ws.Range('E27').Value = '  -1.74%  '

$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '24.79'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +4.51%  '

$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '36.75'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  -3.92%  '

$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '9.62'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +0.82%  '

$ws.Range('E31').Value = '  +0.13%  '

$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '167.83'
$c.Style = 'Normal'

$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '5.29'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +0.44%  '

$ws.Range('E34').Value = '  +0.06%  '

$ws.Range('E35').Value = '  +6.86%  '

$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '3.10'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -2.90%  '

$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '0.0744'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  +0.54%  '

$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '17.52'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -0.41%  '

$ws.Range('E39').Value = '  +1.59%  '

$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '1.87'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +2.17%  '

$ws.Range('E41').Value = '  +1.57%  '

$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '4.26'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +4.05%  '

$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '2.31'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -0.82%  '

$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '19.50'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +2.48%  '

$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.0290'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +2.28%  '

$ws.Range('D46').Value = '1.970.99'
$ws.Range('E46').Value = '  +1.00%  '

$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '2.98'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +2.40%  '

$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '9.95'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -0.18%  '

$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '55.49'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  +3.12%  '

$ws.Range('E50').Value = '  +5.87%  '

$ws.Range('E51').Value = '  +6.76%  '
